$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.452.58'
$ws.Range("E2").Value = '  -2.72%  '
$ws.Range("D3").Value = '''1.775.15'
$ws.Range("E3").Value = '  -1.70%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").Value = '''306.98'
$ws.Range("E6").Value = '  -1.03%  '
$ws.Range("E7").Value = '  +1.63%  '
$ws.Range("E8").Value = '  +2.23%  '
$ws.Range("D9").Value = '''0.07154'
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("D10").Value = '''0.8404'
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("D11").Value = '''20.49'
$ws.Range("E11").Value = '  +1.80%  '
$ws.Range("D12").Value = '''1.804.28'
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = '''6.448'
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("D14").Value = '''5.246'
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").Value = '''0.06892'
$ws.Range("E15").Value = '  +0.12%  '
$ws.Range("D16").Value = '''1.008'
$ws.Range("E16").Value = '  +0.29%  '
$ws.Range("D17").Value = '''78.81'
$ws.Range("E17").Value = '  -3.07%  '
$ws.Range("D18").Value = '''0.000008699'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").Value = '''1.003'
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").Value = '''14.94'
$ws.Range("E20").Value = '  -0.69%  '
$ws.Range("D21").Value = '''26.461.04'
$ws.Range("E21").Value = '  -2.62%  '
$ws.Range("D22").Value = '''5.109'
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("D23").Value = '''11.07'
$ws.Range("E23").Value = '  +2.17%  '
$ws.Range("D24").Value = '''2.000.64'
$ws.Range("E24").Value = '  -1.56%  '
$ws.Range("D25").Value = '''152.56'
$ws.Range("E25").Value = '  -0.64%  '
$ws.Range("D26").Value = '''1.858'
$ws.Range("E26").Value = '  -5.12%  '
$ws.Range("D27").Value = '''18.00'
$ws.Range("E27").Value = '  -1.10%  '
$ws.Range("D28").Value = '''5.062'
$ws.Range("E28").Value = '  +0.47%  '
$ws.Range("D29").Value = '''113.84'
$ws.Range("E29").Value = '  +0.61%  '
$ws.Range("D30").Value = '''1.777'
$ws.Range("E30").Value = '  +4.05%  '
$ws.Range("D31").Value = '''0.08896'
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("D32").Value = '''0.7235'
$ws.Range("E32").Value = '  -1.95%  '
$ws.Range("E33").Value = '  +1.24%  '
$ws.Range("D34").Value = '''4.312'
$ws.Range("E34").Value = '  -3.20%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.751'
$ws.Range("E35").Value = '  -5.51%  '
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = '''1.004'
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").Value = '''1.098'
$ws.Range("E37").Value = '  +2.99%  '
$ws.Range("D38").Value = '''0.05131'
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("D39").Value = '''0.01889'
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").Value = '''0.1612'
$ws.Range("E40").Value = '  -1.34%  '
$ws.Range("D41").Value = '''0.4909'
$ws.Range("E41").Value = '  -0.96%  '
$ws.Range("D42").Value = '''2.593'
$ws.Range("E42").Value = '  -4.13%  '
$ws.Range("D43").Value = '''6.327'
$ws.Range("E43").Value = '  +1.05%  '
$ws.Range("D44").Value = '''7.963'
$ws.Range("E44").Value = '  -2.43%  '
$ws.Range("D45").Value = '''104.68'
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").Value = '''1.005'
$ws.Range("E46").Value = '  +0.33%  '
$ws.Range("D47").Value = '''10.10'
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("D48").Value = '''1.626'
$ws.Range("E48").Value = '  +2.27%  '
$ws.Range("D49").Value = '''0.06178'
$ws.Range("E49").Value = '  -3.21%  '
$ws.Range("D50").Value = '''0.4466'
$ws.Range("E50").Value = '  -1.79%  '
$ws.Range("D51").Value = '''1.712'
$ws.Range("E51").Value = '  +2.07%  '
